$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('AB2').Value = 'maa://21246 (91.26), maa://36684 (98.72), ***maa://22731 (6.67)'
$ws.Range('AF2').Value = 'maa://25251 (92.77), ***maa://21730 (16.92), ***maa://39501 (18.75), *maa://36675 (60.0)'
$ws.Range('D3').Value = 'maa://36987 (93.62), maa://40192 (96.0), maa://39849 (88.89)'
$ws.Range('L3').Value = '*maa://22880 (69.57), maa://20276 (83.56), *maa://22749 (66.67)'
$ws.Range('P3').Value = 'maa://21249 (94.79), maa://26254 (95.65)'
$ws.Range('T3').Value = 'maa://24617 (88.46), **maa://20790 (43.94), ***maa://37170 (19.57)'
$ws.Range('X3').Value = 'maa://27396 (85.52), maa://27484 (95.79), maa://27480 (82.35)'
$ws.Range('D4').Value = 'maa://24632 (93.53), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (97.73), maa://22754 (91.67), maa://27295 (82.14), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('X4').Value = '**maa://32495 (47.01), ***maa://31785 (21.55), ***maa://36683 (28.26)'
$ws.Range('AF4').Value = '*maa://30062 (61.36), ***maa://26209 (13.04), *maa://39394 (76.47)'
$ws.Range('D5').Value = 'maa://21245 (82.32), maa://22744 (83.33)'
$ws.Range('L5').Value = '*maa://22757 (75.86)'
$ws.Range('P5').Value = 'maa://21919 (95.83), maa://21281 (92.31)'
$ws.Range('D6').Value = 'maa://42407 (89.47)'
$ws.Range('X7').Value = 'maa://22399 (94.78), *maa://22758 (71.43)'
$ws.Range('A8').Value = '更新日期：2024.11.03 14:53:27'
$ws.Range('L9').Value = 'maa://22762 (91.57), maa://39552 (88.89)'
$ws.Range('X9').Value = 'maa://26223 (97.14)'
$ws.Range('AB9').Value = 'maa://28711 (88.51), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), **maa://39938 (42.11), maa://40166 (85.71)'
$ws.Range('AF9').Value = 'maa://26206 (90.32), **maa://22865 (47.92)'
$ws.Range('D10').Value = '***maa://25695 (18.86), **maa://32237 (37.84), ***maa://34206 (18.18), ***maa://39951 (16.13), **maa://39243 (33.33)'
$ws.Range('P10').Value = 'maa://28977 (93.59), maa://36669 (85.19), *maa://23264 (61.82)'
$ws.Range('T11').Value = 'maa://22747 (93.75), maa://22501 (98.15)'
$ws.Range('X12').Value = 'maa://22753 (91.5), *maa://21485 (77.1), maa://37962 (84.21)'
$ws.Range('AF12').Value = '*maa://28932 (77.69), *maa://20106 (63.64), *maa://22769 (64.29)'
$ws.Range('P13').Value = 'maa://22676 (92.08), *maa://22583 (75.41), *maa://22500 (56.82)'
$ws.Range('AF13').Value = '**maa://22737 (30.6), maa://39883 (91.43), *maa://39885 (66.67)'
$ws.Range('P14').Value = 'maa://23250 (98.54), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('D15').Value = '*maa://22743 (77.35), maa://22734 (83.48), *maa://30808 (65.0), ***maa://36048 (12.12)'
$ws.Range('H15').Value = 'maa://24304 (88.36), maa://21478 (91.18)'
$ws.Range('AF15').Value = 'maa://21364 (80.68), *maa://22766 (71.43), *maa://36666 (78.26)'
$ws.Range('D16').Value = 'maa://21441 (96.17), maa://36679 (91.43), maa://37650 (96.43)'
$ws.Range('AF16').Value = '*maa://23911 (63.16), maa://27755 (92.21)'
$ws.Range('H17').Value = 'maa://22430 (88.2), maa://39599 (85.19)'
$ws.Range('L17').Value = '*maa://21679 (76.0)'
$ws.Range('T17').Value = '***maa://42324 (26.67)'
$ws.Range('T18').Value = 'maa://24385 (96.97)'
$ws.Range('L20').Value = 'maa://41331 (83.61)'
$ws.Range('T20').Value = 'maa://29113 (88.0)'
$ws.Range('AB21').Value = '*maa://21443 (79.1), ***maa://23820 (29.82)'
$ws.Range('L22').Value = 'maa://27127 (85.87), *maa://22751 (76.19)'
$ws.Range('L23').Value = 'maa://39756 (93.14), maa://39875 (94.34)'
$ws.Range('P23').Value = 'maa://30587 (91.33), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (77.27)'
$ws.Range('X24').Value = 'maa://29988 (85.78), maa://23504 (92.76), **maa://22892 (40.14), *maa://25141 (77.24), *maa://36663 (79.66), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (94.92)'
$ws.Range('H25').Value = '*maa://29063 (74.64), *maa://25311 (74.74), ***maa://22725 (4.84)'
$ws.Range('AF25').Value = 'maa://20108 (96.18), maa://24621 (96.46), maa://36676 (100.0), maa://22771 (85.71), maa://37772 (100.0)'
$ws.Range('AB26').Value = 'maa://42235 (88.1)'
$ws.Range('D28').Value = 'maa://24465 (90.5), maa://25725 (82.93)'
$ws.Range('X28').Value = 'maa://39929 (87.76), ***maa://39723 (14.29), maa://41749 (83.33)'
$ws.Range('AF28').Value = 'maa://36660 (92.73), *maa://36701 (62.96)'
$ws.Range('L29').Value = 'maa://28432 (93.62), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range('AF29').Value = '*maa://24080 (69.17), ***maa://34960 (8.7), maa://42865 (90.0)'
$ws.Range('AB30').Value = 'maa://42979 (92.86)'
$ws.Range('L31').Value = 'maa://35926 (93.39), *maa://36258 (79.75)'
$ws.Range('T32').Value = 'maa://41108 (88.89), maa://41238 (95.24), maa://42859 (94.12)'
$ws.Range('L35').Value = 'maa://41296 (98.63)'
$ws.Range('H41').Value = 'maa://24466 (92.86)'
$ws.Range('H43').Value = 'maa://22525 (92.86), maa://21284 (82.93)'
$ws.Range('H46').Value = 'maa://35931 (92.18)'
$ws.Range('H47').Value = 'maa://27410 (95.89), maa://29661 (97.73), maa://28038 (84.62)'
$ws.Range('H55').Value = 'maa://32532 (92.27)'
$ws.Range('H57').Value = 'maa://25176 (98.0)'
